$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-11 (A:T) with the refreshed TPM-based NATMI values and the
# corrected sending/target cluster labels.
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2087556666666667
$ws.Range("H2").Value = 0.626267
$ws.Range("I2").Value = 0.004212716022507852
$ws.Range("J2").Value = 0.004233737959305679
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.375733333333333
$ws.Range("N2").Value = 4.1272
$ws.Range("O2").Value = 0.457732955319909
$ws.Range("P2").Value = 0.457732955319909
$ws.Range("Q2").Value = 0.2871921291555556
$ws.Range("R2").Value = 2.5847291624
$ws.Range("S2").Value = 0.001928298954906052
$ws.Range("T2").Value = 0.001937921388163069
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2087556666666667
$ws.Range("H3").Value = 0.626267
$ws.Range("I3").Value = 0.004212716022507852
$ws.Range("J3").Value = 0.004233737959305679
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.629803666666667
$ws.Range("N3").Value = 4.889411000000001
$ws.Range("O3").Value = 0.542267044680091
$ws.Range("P3").Value = 0.542267044680091
$ws.Range("Q3").Value = 0.3402307509707778
$ws.Range("R3").Value = 3.062076758737001
$ws.Range("S3").Value = 0.002284417067601801
$ws.Range("T3").Value = 0.00229581657114261
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 48.06441766666666
$ws.Range("H4").Value = 144.193253
$ws.Range("I4").Value = 0.9699460888896084
$ws.Range("J4").Value = 0.9747862312749473
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.375733333333333
$ws.Range("N4").Value = 4.1272
$ws.Range("O4").Value = 0.457732955319909
$ws.Range("P4").Value = 0.457732955319909
$ws.Range("Q4").Value = 66.1238215312889
$ws.Range("R4").Value = 595.1143937816
$ws.Range("S4").Value = 0.4439762897684276
$ws.Range("T4").Value = 0.4461917824466379
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 48.06441766666666
$ws.Range("H5").Value = 144.193253
$ws.Range("I5").Value = 0.9699460888896084
$ws.Range("J5").Value = 0.9747862312749473
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.629803666666667
$ws.Range("N5").Value = 4.889411000000001
$ws.Range("O5").Value = 0.542267044680091
$ws.Range("P5").Value = 0.542267044680091
$ws.Range("Q5").Value = 78.33556414933145
$ws.Range("R5").Value = 705.0200773439831
$ws.Range("S5").Value = 0.5259697991211808
$ws.Range("T5").Value = 0.5285944488283094
$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.245141
$ws.Range("H6").Value = 0.735423
$ws.Range("I6").Value = 0.004946976697512072
$ws.Range("J6").Value = 0.004971662679410635
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.375733333333333
$ws.Range("N6").Value = 4.1272
$ws.Range("O6").Value = 0.457732955319909
$ws.Range("P6").Value = 0.457732955319909
$ws.Range("Q6").Value = 0.3372486450666667
$ws.Range("R6").Value = 3.0352378056
$ws.Range("S6").Value = 0.002264394263650924
$ws.Range("T6").Value = 0.002275693851100327
$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Cxcl13"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.245141
$ws.Range("H7").Value = 0.735423
$ws.Range("I7").Value = 0.004946976697512072
$ws.Range("J7").Value = 0.004971662679410635
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.629803666666667
$ws.Range("N7").Value = 4.889411000000001
$ws.Range("O7").Value = 0.542267044680091
$ws.Range("P7").Value = 0.542267044680091
$ws.Range("Q7").Value = 0.3995317006503334
$ws.Range("R7").Value = 3.595785305853001
$ws.Range("S7").Value = 0.002682582433861148
$ws.Range("T7").Value = 0.002695968828310308
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cxcl13"
$ws.Range("C8").Value = "Cxcr3"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.7381525
$ws.Range("H8").Value = 1.476305
$ws.Range("I8").Value = 0.01489601175123818
$ws.Range("J8").Value = 0.00998022970715808
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.375733333333333
$ws.Range("N8").Value = 4.1272
$ws.Range("O8").Value = 0.457732955319909
$ws.Range("P8").Value = 0.457732955319909
$ws.Range("Q8").Value = 1.015500999333333
$ws.Range("R8").Value = 6.093005996
$ws.Range("S8").Value = 0.006818395481374345
$ws.Range("T8").Value = 0.004568280038629018
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cxcl13"
$ws.Range("C9").Value = "Cxcr3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.7381525
$ws.Range("H9").Value = 1.476305
$ws.Range("I9").Value = 0.01489601175123818
$ws.Range("J9").Value = 0.00998022970715808
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.629803666666667
$ws.Range("N9").Value = 4.889411000000001
$ws.Range("O9").Value = 0.542267044680091
$ws.Range("P9").Value = 0.542267044680091
$ws.Range("Q9").Value = 1.203043651059167
$ws.Range("R9").Value = 7.218261906355001
$ws.Range("S9").Value = 0.008077616269863835
$ws.Range("T9").Value = 0.005411949668529063
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Cxcl13"
$ws.Range("C10").Value = "Cxcr3"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2972333333333333
$ws.Range("H10").Value = 0.8917
$ws.Range("I10").Value = 0.005998206639133552
$ws.Range("J10").Value = 0.006028138379178328
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.375733333333333
$ws.Range("N10").Value = 4.1272
$ws.Range("O10").Value = 0.457732955319909
$ws.Range("P10").Value = 0.457732955319909
$ws.Range("Q10").Value = 0.4089138044444445
$ws.Range("R10").Value = 3.68022424
$ws.Range("S10").Value = 0.0027455768515501
$ws.Range("T10").Value = 0.002759277595378662
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Cxcl13"
$ws.Range("C11").Value = "Cxcr3"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.2972333333333333
$ws.Range("H11").Value = 0.8917
$ws.Range("I11").Value = 0.005998206639133552
$ws.Range("J11").Value = 0.006028138379178328
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.629803666666667
$ws.Range("N11").Value = 4.889411000000001
$ws.Range("O11").Value = 0.542267044680091
$ws.Range("P11").Value = 0.542267044680091
$ws.Range("Q11").Value = 0.4844319765222223
$ws.Range("R11").Value = 4.359887788700001
$ws.Range("S11").Value = 0.003252629787583452
$ws.Range("T11").Value = 0.003268860783799666


# Rows 12-13 no longer exist in the refreshed export (Resolving-Mac as
# sender is now fully covered by rows 10-11); remove them so the sheet
# shrinks from A1:T13 to A1:T11.
$ws.Range("A12:T13").EntireRow.Delete()
